$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 8.720000000000001
$ws.Range("E2").Value = 59.3
$ws.Range("F2").Value = 11.87
$ws.Range("N2").Value = 50.68470204858703

$ws.Range("D3").Value = 15.38
$ws.Range("E3").Value = 57.5
$ws.Range("F3").Value = 6.58
$ws.Range("N3").Value = 50.68470204858703
